{"js": "// Translate the English subtitle lines into Swahili (Kenya).\n// Each entry is an exact, whole-paragraph text replacement.\nconst replacements = [\n  [\"The playful mathematicians:\", \"Wanahisabati wanaocheza:\"],\n  [\n    \"** the dialogue starts at second 47, so I added 28 seconds to all the times as they were. -John Argentino\",\n    \"** mazungumzo huanza saa 47 ya pili, kwa hivyo niliongeza sekunde 28 kwa nyakati zote kama zilivyokuwa. -John Argentino\",\n  ],\n  [\"there are two mathematicians, let's call\", \"kuna wanahisabati wawili, tupige simu\"],\n  [\"them Fil and Mike who meet each other\", \"Fil na Mike wanaokutana\"],\n  [\"again after a long time. After some\", \"tena baada ya muda mrefu. Baada ya baadhi\"],\n  [\"chatting, Phil says he has three children, then\", \"kuzungumza, Phil anasema ana watoto watatu, basi\"],\n  [\"Mike, astonished, asks: 'How old are they?' Fil,\", \"Kwa mshangao, Mike anauliza: 'Wana umri gani?' Fil,\"],\n  [\"being a playful mathematician, answers\", \"kuwa mwanahisabati mchezaji, anajibu\"],\n  [\"'You tell me! I'll give you a hint: if you\", \"'Wewe niambie! Nitakupa kidokezo: ikiwa wewe\"],\n  [\"multiply the three ages together you\", \"zidisheni enzi tatu pamoja ninyi\"],\n  [\"get 36.' Mike takes sometimes to think\", \"pata 36.' Mike huchukua wakati mwingine kufikiria\"],\n  [\"and says: 'I'm sorry Fil, but I do need\", \"na kusema: 'Samahani Fil, lakini nahitaji\"],\n  [\"another hint. So Fil tells Mike:\", \"kidokezo kingine. Kwa hivyo Fil anamwambia Mike:\"],\n  [\"'Yes, sure, here it is: if you had up to\", \"'Ndiyo, hakika, hapa ni: kama alikuwa na hadi\"],\n  [\"three ages you get the number of math\", \"miaka mitatu unapata idadi ya hesabu\"],\n  [\"papers we publish together. Do you remember it?'\", \"karatasi tunachapisha pamoja. Je, unaikumbuka?'\"],\n  [\"'Yes I do remember How many, but still\", \"'Ndio nakumbuka wangapi, lakini bado\"],\n  [\"I do not have enough information! I need\", \"Sina taarifa za kutosha! nahitaji\"],\n  [\"at least one more.' Fil says: 'Yes don't\", \"angalau moja zaidi.' Fil anasema: 'Ndiyo usifanye hivyo\"],\n  [\"worry but this is the last one:\", \"wasiwasi lakini hii ni ya mwisho:\"],\n  [\"The youngest one has blues eyes.' And\", \"Mdogo ana macho ya blues.' Na\"],\n  [\"suddenly Mike gets the answer. You\", \"ghafla Mike anapata jibu. Wewe\"],\n  [\"hear the conversation but you don't know\", \"sikia mazungumzo lakini hujui\"],\n  [\"how many papers they published together.\", \"ni karatasi ngapi walichapisha pamoja.\"],\n  [\"However, you do want to know the ages of\", \"Hata hivyo, unataka kujua umri wa\"],\n  [\"the three children. Can you figure them\", \"watoto watatu. Je, unaweza kuwahesabu\"],\n  [\"out?\", \"nje?\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text;\n  for (const [oldText, newText] of replacements) {\n    if (text === oldText) {\n      paragraph.getRange().insertText(newText, \"Replace\");\n      break;\n    }\n  }\n}\nawait context.sync();\n\n// \"[Music]\" occurs twice, and both must become \"[Muziki]\".\nconst musicResults = body.search(\"[Music]\", { matchCase: true });\nmusicResults.load(\"items\");\nawait context.sync();\nfor (const range of musicResults.items) {\n  range.insertText(\"[Muziki]\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Translate the English subtitle lines into Swahili (Kenya).\n# Each entry is an exact, whole-paragraph text replacement: Find() locates\n# the English text and the matched Range's .Text is then overwritten with\n# the Swahili translation.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"The playful mathematicians:\", \"Wanahisabati wanaocheza:\"),\n    @(\"** the dialogue starts at second 47, so I added 28 seconds to all the times as they were. -John Argentino\", \"** mazungumzo huanza saa 47 ya pili, kwa hivyo niliongeza sekunde 28 kwa nyakati zote kama zilivyokuwa. -John Argentino\"),\n    @(\"[Music]\", \"[Muziki]\"),\n    @(\"there are two mathematicians, let's call\", \"kuna wanahisabati wawili, tupige simu\"),\n    @(\"them Fil and Mike who meet each other\", \"Fil na Mike wanaokutana\"),\n    @(\"again after a long time. After some\", \"tena baada ya muda mrefu. Baada ya baadhi\"),\n    @(\"chatting, Phil says he has three children, then\", \"kuzungumza, Phil anasema ana watoto watatu, basi\"),\n    @(\"Mike, astonished, asks: 'How old are they?' Fil,\", \"Kwa mshangao, Mike anauliza: 'Wana umri gani?' Fil,\"),\n    @(\"being a playful mathematician, answers\", \"kuwa mwanahisabati mchezaji, anajibu\"),\n    @(\"'You tell me! I'll give you a hint: if you\", \"'Wewe niambie! Nitakupa kidokezo: ikiwa wewe\"),\n    @(\"multiply the three ages together you\", \"zidisheni enzi tatu pamoja ninyi\"),\n    @(\"get 36.' Mike takes sometimes to think\", \"pata 36.' Mike huchukua wakati mwingine kufikiria\"),\n    @(\"and says: 'I'm sorry Fil, but I do need\", \"na kusema: 'Samahani Fil, lakini nahitaji\"),\n    @(\"another hint. So Fil tells Mike:\", \"kidokezo kingine. Kwa hivyo Fil anamwambia Mike:\"),\n    @(\"'Yes, sure, here it is: if you had up to\", \"'Ndiyo, hakika, hapa ni: kama alikuwa na hadi\"),\n    @(\"three ages you get the number of math\", \"miaka mitatu unapata idadi ya hesabu\"),\n    @(\"papers we publish together. Do you remember it?'\", \"karatasi tunachapisha pamoja. Je, unaikumbuka?'\"),\n    @(\"'Yes I do remember How many, but still\", \"'Ndio nakumbuka wangapi, lakini bado\"),\n    @(\"I do not have enough information! I need\", \"Sina taarifa za kutosha! nahitaji\"),\n    @(\"at least one more.' Fil says: 'Yes don't\", \"angalau moja zaidi.' Fil anasema: 'Ndiyo usifanye hivyo\"),\n    @(\"worry but this is the last one:\", \"wasiwasi lakini hii ni ya mwisho:\"),\n    @(\"The youngest one has blues eyes.' And\", \"Mdogo ana macho ya blues.' Na\"),\n    @(\"suddenly Mike gets the answer. You\", \"ghafla Mike anapata jibu. Wewe\"),\n    @(\"hear the conversation but you don't know\", \"sikia mazungumzo lakini hujui\"),\n    @(\"how many papers they published together.\", \"ni karatasi ngapi walichapisha pamoja.\"),\n    @(\"However, you do want to know the ages of\", \"Hata hivyo, unataka kujua umri wa\"),\n    @(\"the three children. Can you figure them\", \"watoto watatu. Je, unaweza kuwahesabu\"),\n    @(\"out?\", \"nje?\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    # Locate each occurrence with Find (no Replacement text is supplied,\n    # so Word's \"smart quotes\" AutoCorrect can't silently rewrite straight\n    # apostrophes in $replaceText into curly quotes), then overwrite the\n    # matched range's Text directly - a plain, literal text substitution.\n    # \"[Music]\" matches twice, everything else matches exactly once, so\n    # loop over every hit found while scanning forward from the doc start.\n    $searchRange = $d.Content\n    $find = $searchRange.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 0\n\n    while ($find.Execute()) {\n        $searchRange.Text = $replaceText\n        $searchRange.SetRange($searchRange.End, $d.Content.End)\n    }\n}\n"}
